# CIERRE DEL 23 OCT 2021
# Advance the weekly payroll sheet from SEMANA 42 (Oct 11-17, 2021) to
# SEMANA 43 (Oct 18-24, 2021) and update this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Week banner --------------------------------------------------------
# B9 holds the shared-string week banner; H9, B27, H27, B43, H43, B60 all
# pull from it via formulas (=B9, =B27, =H27, =B43, =H43) and recalc
# automatically.
$ws.Range("B9").Value = "SEMANA   43  DEL    18      Al    24   DE   OCTUBRE          2021"

# --- Employee #1 (JUSTINA REYES LEAL) -----------------------------------
# "Extra" pay goes to 0 this week; K6 = SUM(K3:K5) recalcs automatically.
$ws.Range("K4").Value = 0

# --- Employee #3 (CLAUDIA PRIETO VARGAS) --------------------------------
# Days pay bumped to 2800; K24 = SUM(K21:K23) recalcs automatically.
$ws.Range("K21").Value = 2800

# --- Employee #6 (Ma. Margarita AGUILAR AQUINO) -------------------------
# Prestamo (loan) deduction cleared, Extras paid out; K41 = SUM(K38:K40)
# recalcs automatically.
$ws.Range("K39").Value = 0
$ws.Range("K40").Value = 1250

# --- Window / selection state -------------------------------------------
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("H60:H61").Select()
